$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Vreme skeniranja" (scan time) timestamps in column E
$ws.Range("E2").Value = "2025-04-14 19:42:47"
$ws.Range("E3").Value = "2025-04-14 19:42:46"
$ws.Range("E4").Value = "2025-04-14 19:42:47"
$ws.Range("E5").Value = "2025-04-14 19:42:49"
$ws.Range("E6").Value = "2025-04-14 19:42:51"
$ws.Range("E7").Value = "2025-04-14 19:42:53"
$ws.Range("E8").Value = "2025-04-14 19:42:54"
$ws.Range("E9").Value = "2025-04-14 19:42:55"
$ws.Range("E10").Value = "2025-04-14 19:42:55"
$ws.Range("E11").Value = "2025-04-14 19:42:57"
$ws.Range("E12").Value = "2025-04-14 19:42:58"
$ws.Range("E13").Value = "2025-04-14 19:42:59"

# Fix "Kutija" (box) id for row 13
$ws.Range("B13").Value = 26002680638
